# ---------------------------------------------------------------------------
# Bugfixed the naive forecaster component module
#
# The underlying AR(2) YoY naive-forecaster pipeline was re-run with the
# corrected rolling window: every (date_of_forecast, y_0, y_0_forecast,
# y_1, y_1_forecast) observation shifts one row earlier (the oldest
# observation drops off the front of the window) and the forecast columns
# (y_0_forecast / y_1_forecast) are recomputed from the fixed component.
# This refreshes the full data block and removes the now out-of-window
# trailing row from the sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The corrected series only spans 51 data rows (plus header) instead of 52,
# so drop the last (now stale) row of the old series first.
$ws.Rows("53").Delete()

# Refreshed values for A2:E52 -> (date_of_forecast, y_0, y_0_forecast, y_1, y_1_forecast)
$rows = @(
    @(39583, 2008, $null, 2009, $null),
    @(39765, 2008, $null, 2009, $null),
    @(39948, 2009, $null, 2010, $null),
    @(40130, 2009, -0.5555135891318952, 2010, $null),
    @(40310, 2010, $null, 2011, $null),
    @(40494, 2010, 0.8442071301477228, 2011, $null),
    @(40676, 2011, $null, 2012, $null),
    @(40862, 2011, 1.122475521884692, 2012, $null),
    @(41044, 2012, $null, 2013, $null),
    @(41228, 2012, -0.578174579726376, 2013, -0.3496173419443749),
    @(41409, 2013, -0.6616365666142765, 2014, 0.2197847717222867),
    @(41592, 2013, -0.7492845378401558, 2014, 0.05500386022236903),
    @(41774, 2014, 0.2184978785563896, 2015, 0.04195831742983547),
    @(41957, 2014, 0.2751437421933511, 2015, 0.07916875696107883),
    @(42137, 2015, -0.01790997771649039, 2016, 0.1671491311400208),
    @(42321, 2015, 0.07468705617190707, 2016, 0.1656566557188155),
    @(42503, 2016, -0.08711135105702317, 2017, 0.1158714888162216),
    @(42689, 2016, -0.05493014849097255, 2017, 0.2044493994367125),
    @(42867, 2017, 0.1761917659537371, 2018, 0.1078587431702305),
    @(43053, 2017, 0.2820931576894115, 2018, 0.2685680645708288),
    @(43145, 2018, 0.2803378563356329, 2019, 0.1089201880626334),
    @(43235, 2018, 0.3524405906205841, 2019, 0.1531204771924033),
    @(43326, 2018, 0.1415392254179304, 2019, -0.1904239862803969),
    @(43418, 2018, 0.2343541283920114, 2019, 0.2977174885593792),
    @(43510, 2019, 0.3424613118119479, 2020, 0.1339497680586277),
    @(43600, 2019, 0.07560805834034845, 2020, -0.01670081902098719),
    @(43691, 2019, -0.05069288950212414, 2020, -0.2043633904885378),
    @(43783, 2019, -0.009430310228020211, 2020, 0.0103609600907939),
    @(43875, 2020, 0.01105513701109562, 2021, 0.0925667197466451),
    @(43966, 2020, -0.5849047489490333, 2021, -0.2059746096811033),
    @(44068, 2020, -2.657403949513992, 2021, -1.40802832891157),
    @(44159, 2020, -2.657403949513992, 2021, -1.407243743159736),
    @(44251, 2021, -0.2500618974080826, 2022, 0.07829984441984905),
    @(44341, 2021, -0.4334047671505248, 2022, -0.06335028919957075),
    @(44432, 2021, -0.3096364143617802, 2022, 0.2283024244226883),
    @(44525, 2021, -0.3096364143617802, 2022, 0.2048390592685578),
    @(44617, 2022, -0.0355780787674953, 2023, -0.1890623092888566),
    @(44706, 2022, -0.1663214453978101, 2023, -0.3613518455741316),
    @(44798, 2022, -0.1730430455425092, 2023, -0.4115424244148125),
    @(44890, 2022, -0.1730430455425092, 2023, 0.2152263639657814),
    @(44981, 2023, 0.6050248749486009, 2024, -0.1865259660156937),
    @(45071, 2023, 0.6502606143725664, 2024, -0.1691853834640433),
    @(45163, 2023, 0.6376744206510576, 2024, -0.2058599286704377),
    @(45254, 2023, 0.6376744206510576, 2024, -0.06071040501895997),
    @(45345, 2024, -0.05036452040672046, 2025, -0.2233431232791294),
    @(45436, 2024, 0.1549171986535924, 2025, -0.03486668218654065),
    @(45534, 2024, 0.1856341247700399, 2025, 0.0444032571666142),
    @(45618, 2024, 0.1856341247700399, 2025, 0.160714157635633),
    @(45713, 2025, 0.01130372647704103, 2026, -0.1505299872523014),
    @(45800, 2025, -0.06391119588061711, 2026, -0.2122873162357264),
    @(45891, 2025, -0.09450306168263811, 2026, -0.2824524929558314)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
